# Updates cryptocurrency price (D) and 1h-volume-change (E) figures in the
# cryptos worksheet (Sat Dec 30 22:27:48 UTC 2023 GitHub Actions refresh).
#
# Price cells are stored as plain text in the workbook (some values use a
# "." thousands separator, e.g. "42.451.73", which would be misread as a
# number). For values that would otherwise look like a plain decimal number
# (e.g. "316.33"), the cell's NumberFormat is forced to "@" (Text) right
# before the value is written so it round-trips as a text cell instead of
# being coerced into a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.451.73"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "2.298.66"
$ws.Range("E3").Value = "  +1.24%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.33"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.64"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.42"
$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.960"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.25"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("D16").Value = "2.648.83"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "2.300.14"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").Value = "42.418.21"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.44"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.82"
$ws.Range("E22").Value = "  +28.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  +3.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "276.02"
$ws.Range("E24").Value = "  +8.46%  "

$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.82"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.75"
$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.41"
$ws.Range("E30").Value = "  +6.18%  "

$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0875"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.92"
$ws.Range("E33").Value = "  +3.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("E34").Value = "  +4.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.119"
$ws.Range("E35").Value = "  +2.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.63"
$ws.Range("E36").Value = "  -9.92%  "

$ws.Range("E37").Value = "  +1.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0364"
$ws.Range("E38").Value = "  +4.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  +3.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("E41").Value = "  +3.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.05"
$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "94.62"
$ws.Range("E43").Value = "  -2.14%  "

$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.07"
$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "80.80"
$ws.Range("E47").Value = "  +9.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.44"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("D51").Value = "1.587.21"
$ws.Range("E51").Value = "  +2.05%  "
